# Fixed errors in tables
# This script reproduces the commit "Fixed errors in tables":
#  - swaps the Condition/Correct columns (C/D) in both result tables
#    (including the header row and a handful of data rows that had been
#    entered in the wrong column)
#  - fixes a couple of mis-keyed data points
#  - repairs the SUM/AVERAGE helper formulas that had drifted after the
#    column swap
#  - re-applies the correct cell styling (borders/fills/alignment) that
#    Excel had renumbered when it resaved the workbook
#  - restores the saved selection / scroll position

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Swap the "Condition" / "Correct" column contents (columns C and D)
#    for every row where the two values need to trade places.
# ---------------------------------------------------------------------
function Swap-Values($ref1, $ref2) {
    $r1 = $ws.Range($ref1)
    $r2 = $ws.Range($ref2)
    $tmp = $r1.Value2
    $r1.Value2 = $r2.Value2
    $r2.Value2 = $tmp
}

# Header rows (text swap: "Condition" <-> "Correct")
Swap-Values "C4" "D4"
Swap-Values "C21" "D21"

# Data rows that had Condition/Correct entered in the wrong column
Swap-Values "C5"  "D5"
Swap-Values "C8"  "D8"
Swap-Values "C9"  "D9"
Swap-Values "C13" "D13"
Swap-Values "C22" "D22"
Swap-Values "C28" "D28"
Swap-Values "C31" "D31"

# ---------------------------------------------------------------------
# 2. Correct a couple of mis-entered data points
# ---------------------------------------------------------------------
$ws.Range("G26").Value2 = 0
$ws.Range("E30").Value2 = 1

# ---------------------------------------------------------------------
# 3. Repair the helper formulas in the first table (rows 15-17) so the
#    totals/averages reference the corrected columns.
# ---------------------------------------------------------------------
$ws.Range("C15").Formula = "=SUM(C5:C14)"
$ws.Range("D15").Formula = "=SUM(D5:D14)"
$ws.Range("E15").Formula = "=SUM(E5:E14)"
$ws.Range("F15").Formula = "=SUM(F5:F14)"
$ws.Range("G15").Formula = "=SUM(G5:G14)"

$ws.Range("C16").Formula = "=SUM(D15:G15)"
$ws.Range("C17").Formula = "=AVERAGE(D5:G14)"

# ---------------------------------------------------------------------
# 4. Re-apply cell styling (borders / fills / alignment) for the cells
#    whose style index changed when Excel consolidated the style table.
# ---------------------------------------------------------------------

# Colors measured from the workbook's existing theme fills
$FillBlue   = 11854022   # theme 4 (accent1), tint 0.6  -- header band fill
$FillOrange = 15652797   # theme 9 (accent6), tint 0.6  -- section title fill

function Set-ThinBox($range) {
    foreach ($cell in $range.Cells) {
        foreach ($edge in 7, 8, 9, 10) {
            $b = $cell.Borders.Item($edge)
            $b.LineStyle = 1
            $b.Weight = 2
            $b.ColorIndex = 0
        }
    }
}

function Clear-NoFill($range) {
    $range.Interior.Pattern = -4142
    $range.Interior.ColorIndex = -4142
}

function Set-Fill($range, $color) {
    $range.Interior.Pattern = 1
    $range.Interior.Color = $color
}

# B3: plain box border, no fill, general alignment (style 3)
Set-ThinBox $ws.Range("B3")
Clear-NoFill $ws.Range("B3")
$ws.Range("B3").HorizontalAlignment = -4131

# C3: orange fill, centered, border without top (style 6)
Set-ThinBox $ws.Range("C3")
Set-Fill $ws.Range("C3") $FillOrange
$ws.Range("C3").Borders.Item(8).LineStyle = -4142
$ws.Range("C3").HorizontalAlignment = -4108

# D3:G3: orange fill, centered, full box border (style 7)
Set-ThinBox $ws.Range("D3:G3")
Set-Fill $ws.Range("D3:G3") $FillOrange
$ws.Range("D3:G3").HorizontalAlignment = -4108

# B4: box border, no fill, general alignment (style 2)
Set-ThinBox $ws.Range("B4")
Clear-NoFill $ws.Range("B4")
$ws.Range("B4").HorizontalAlignment = -4131

# C4:G4: orange fill, box border, general alignment (style 5)
Set-ThinBox $ws.Range("C4:G4")
Set-Fill $ws.Range("C4:G4") $FillOrange
$ws.Range("C4:G4").HorizontalAlignment = -4131

# C16:G17 and C33:G34: box border, no fill, centered (style 9)
foreach ($addr in "C16:G17", "C33:G34") {
    $rng = $ws.Range($addr)
    Set-ThinBox $rng
    Clear-NoFill $rng
    $rng.HorizontalAlignment = -4108
}

# C20:G20: blue fill, box border, centered (style 8)
Set-ThinBox $ws.Range("C20:G20")
Set-Fill $ws.Range("C20:G20") $FillBlue
$ws.Range("C20:G20").HorizontalAlignment = -4108

# C21:G21: blue fill, box border, general alignment (style 4)
Set-ThinBox $ws.Range("C21:G21")
Set-Fill $ws.Range("C21:G21") $FillBlue
$ws.Range("C21:G21").HorizontalAlignment = -4131

# ---------------------------------------------------------------------
# 5. Restore saved view state (scroll position / selection)
# ---------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("G26").Select() | Out-Null
